$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Entities"
$ws.Activate()
$ws.Range("K16").Select()
